$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (current "Venue" column),
# shifting Venue..geometry from B:M to D:O.
$ws.Range("B:C").Insert()

# New header labels for the inserted columns (match the header style of the
# rest of row 1, which was carried over to column D after the shift).
$ws.Cells.Item(1, 2).Value = "Unnamed: 0.1"
$ws.Cells.Item(1, 3).Value = "Unnamed: 0"
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Fill the new columns with the same row-index values as column A (rows 2-30).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $idxVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $idxVal
    $ws.Cells.Item($r, 3).Value = $idxVal
}

# The column insert caused the new B/C data cells to inherit column A's bold
# header-style formatting; clear it back to the default (unstyled) look used
# by the other data columns.
$ws.Range("D2").Copy()
$ws.Range("B2:C" + $lastRow).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
